$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.732.39'
$ws.Range("E2").Value = '  +4.39%  '

# Row 3
$ws.Range("D3").Value = '2.778.60'
$ws.Range("E3").Value = '  +5.61%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.73'
$ws.Range("E5").Value = '  +3.34%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '333.81'
$ws.Range("E6").Value = '  +3.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.540'
$ws.Range("E7").Value = '  +2.27%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  +6.40%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.07'
$ws.Range("E10").Value = '  +5.87%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0867'
$ws.Range("E11").Value = '  +6.75%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.36'
$ws.Range("E12").Value = '  +2.73%  '

# Row 13
$ws.Range("E13").Value = '  +2.20%  '

# Row 14
$ws.Range("E14").Value = '  +4.96%  '

# Row 15
$ws.Range("D15").Value = '3.215.95'
$ws.Range("E15").Value = '  +5.48%  '

# Row 16
$ws.Range("D16").Value = '2.779.86'
$ws.Range("E16").Value = '  +5.43%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.891'
$ws.Range("E17").Value = '  +3.73%  '

# Row 18
$ws.Range("D18").Value = '51.730.92'
$ws.Range("E18").Value = '  +4.56%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.36'
$ws.Range("E19").Value = '  +13.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.58'
$ws.Range("E20").Value = '  +5.48%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.89'
$ws.Range("E21").Value = '  +3.03%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0979'
$ws.Range("E22").Value = '  +3.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '279.23'
$ws.Range("E23").Value = '  +3.47%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.91'
$ws.Range("E24").Value = '  +1.36%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.70'
$ws.Range("E25").Value = '  +6.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.94'
$ws.Range("E26").Value = '  +2.48%  '

# Row 27
$ws.Range("E27").Value = '  +0.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.21'
$ws.Range("E28").Value = '  -0.90%  '

# Row 29
$ws.Range("E29").Value = '  +0.23%  '

# Row 30
$ws.Range("E30").Value = '  +2.57%  '

# Row 31
$ws.Range("E31").Value = '  +0.91%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.24'
$ws.Range("E32").Value = '  +1.79%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.60'
$ws.Range("E33").Value = '  +1.86%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0823'
$ws.Range("E34").Value = '  +1.21%  '

# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.08'
$ws.Range("E35").Value = '  +3.47%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.17'
$ws.Range("E36").Value = '  +1.19%  '

# Row 37
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
$ws.Range("E38").Value = '  +2.85%  '

# Row 39
$ws.Range("E39").Value = '  +4.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0355'
$ws.Range("E40").Value = '  +8.76%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '127.22'
$ws.Range("E41").Value = '  +0.64%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.36'
$ws.Range("E42").Value = '  +4.15%  '

# Row 43
$ws.Range("E43").Value = '  +3.28%  '

# Row 44
$ws.Range("E44").Value = '  +7.27%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.47'
$ws.Range("E45").Value = '  +16.74%  '

# Row 46
$ws.Range("D46").Value = '2.099.49'
$ws.Range("E46").Value = '  +2.05%  '

# Row 47
$ws.Range("E47").Value = '  +3.41%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.24'
$ws.Range("E48").Value = '  +3.90%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.57'
$ws.Range("E49").Value = '  +6.80%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '60.53'
$ws.Range("E50").Value = '  +2.63%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.89'
$ws.Range("E51").Value = '  -0.39%  '
